# Transcription, JS work, and stylesheet rewrite
# Adds four new glyph entries (g35-g38) to the "Glyphs" sheet, continuing
# the existing A=code / B=description table that runs from row 1 to 35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Glyphs")

# Row 36
$ws.Range("A36").Value = "g35"
$ws.Range("B36").Value = "ligature"

# Row 37 (code filled now, description filled later - out of order,
# matching the shared-string insertion order of the source edit)
$ws.Range("A37").Value = "g36"

# Row 38's description was entered before row 37's description.
$ws.Range("B38").Value = "ac ligature"
$ws.Range("B37").Value = "p with tail"

# Row 38's code, then row 39's code.
$ws.Range("A38").Value = "g37"
$ws.Range("A39").Value = "g38"

# Row 39's description, last.
$ws.Range("B39").Value = "ur abbreviation"

# Leave the sheet scrolled/selected near the newly added rows, matching
# the final cursor position of the source edit (B39).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("B39").Select()
